$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 172, shifting rows 172:191 down to 173:192
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new record
$ws.Cells.Item(172, 1).Value = 3
$ws.Cells.Item(172, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(172, 3).Value = "Coquimbo"
$ws.Cells.Item(172, 4).Value = 44449
$ws.Cells.Item(172, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(172, 5).Value = 5
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100108
$ws.Cells.Item(172, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(172, 9).Value = 100108002
$ws.Cells.Item(172, 10).Value = "Mango"
$ws.Cells.Item(172, 11).Value = "Sin especificar"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 228
$ws.Cells.Item(172, 14).Value = 9000
$ws.Cells.Item(172, 15).Value = 9000
$ws.Cells.Item(172, 16).Value = 9000
$ws.Cells.Item(172, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(172, 18).Value = "Brasil"
$ws.Cells.Item(172, 19).Value = 2250
$ws.Cells.Item(172, 20).Value = 4
